$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 5680
$wsExhibit.Range("F7").Value = 3
$wsExhibit.Range("F8").Value = 54

# Sheet "全部类型" (All types)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 5680
$wsAll.Range("F8").Value = 3
$wsAll.Range("F10").Value = 54
